$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Ingreso")
$ws2 = $wb.Worksheets.Item("Gastos")

# --- Ingreso sheet: add new contribution rows for 2023-07-16 (serial 45123) ---
# Fill rows 483-490 first (order matters for shared-string creation order:
# "chamo" must be created before "Fernando").
$rows483to490 = @(
    @{ Row = 483; Member = "Melvin";     Amount = 100 },
    @{ Row = 484; Member = "Randy";      Amount = 100 },
    @{ Row = 485; Member = "chamo";      Amount = 100 },
    @{ Row = 486; Member = "Invitados";  Amount = 300 },
    @{ Row = 487; Member = "Fernando";   Amount = 100 },
    @{ Row = 488; Member = "Omaury";     Amount = 100 },
    @{ Row = 489; Member = "Anuel";      Amount = 100 },
    @{ Row = 490; Member = "Julio";      Amount = 200 }
)

foreach ($entry in $rows483to490) {
    $r = $entry.Row
    $ws1.Cells.Item($r, 1).Value = 45123
    $ws1.Cells.Item($r, 2).Value = $entry.Member
    $ws1.Cells.Item($r, 3).Value = $entry.Amount
    $ws1.Cells.Item($r, 3).Style = "Normal"
    $ws1.Cells.Item($r, 4).Value = "Aporte"
}

# Now fix the old "Invitados" entry on row 478 to "Fernando" (reuses the
# shared string created above).
$ws1.Cells.Item(478, 2).Value = "Fernando"

# Rows 491-492 keep the inherited "Aporte" numeric style (s="2").
$rows491to492 = @(
    @{ Row = 491; Member = "Jeicol"; Amount = 400 },
    @{ Row = 492; Member = "Robert"; Amount = 250 }
)

foreach ($entry in $rows491to492) {
    $r = $entry.Row
    $ws1.Cells.Item($r, 1).Value = 45123
    $ws1.Cells.Item($r, 2).Value = $entry.Member
    $ws1.Cells.Item($r, 3).Value = $entry.Amount
    $ws1.Cells.Item($r, 4).Value = "Aporte"
}

# --- Gastos sheet: add new expense row for 2023-07-16 ---
$ws2.Cells.Item(51, 1).Value = 45123
$ws2.Cells.Item(51, 2).Value = "Arbitro y agua"
$ws2.Cells.Item(51, 3).Value = 1100

$ws2.Range("A51").Select() | Out-Null

# Leave "Ingreso" as the active/selected sheet and row 492 as the active cell,
# matching the saved view state of the workbook.
$ws1.Activate()
$ws1.Range("A492").Select() | Out-Null
